$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2026-02 (row 27)
$ws.Range("B27").Value = 6550
$ws.Range("C27").Value = 1019
$ws.Range("D27").Value = 6107241
$ws.Range("E27").Value = 932.4032061068702
$ws.Range("F27").Value = 10.08403361344539
$ws.Range("G27").Value = 7.602956705385422
$ws.Range("H27").Value = 25.26371103559075
